$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New team roster (Matteo Mazzola's team) added below the existing rows.
$names = @(
    "Matteo Mazzola",
    "Stefano Pizzini",
    "Marco Gerola",
    "Michele Parisi",
    "Andrea Anzelini",
    "Davide Raffaelli"
)

$startRow = 8
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i

    # Columns A-F stay blank (matching the existing rows' layout) but still
    # need a real (empty-string) cell present, so write an empty text value
    # and then strip the quote-prefix style it introduces.
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = "'"
        $cell.Style = "Normal"
    }

    # Column G holds the player's name.
    $ws.Cells.Item($row, 7).Value = $names[$i]
}
